$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.481.33'
$ws.Range("E2").Value = '  +0.17%  '
$ws.Range("D3").Value = '2.659.24'
$ws.Range("E3").Value = '  +1.16%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''595.81'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").Value = '''175.81'
$ws.Range("E6").Value = '  -2.13%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '''0.522'
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").Value = '2.657.07'
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("E10").Value = '  -2.22%  '
$ws.Range("D11").Value = '''0.169'
$ws.Range("E11").Value = '  +1.77%  '
$ws.Range("E12").Value = '  -0.45%  '
$ws.Range("D13").Value = '''5.00'
$ws.Range("E13").Value = '  -0.66%  '
$ws.Range("D14").Value = '3.145.43'
$ws.Range("E14").Value = '  +0.78%  '
$ws.Range("E15").Value = '  -2.06%  '
$ws.Range("D16").Value = '72.353.34'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '''26.11'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '2.658.36'
$ws.Range("E18").Value = '  +0.37%  '
$ws.Range("D19").Value = '''12.25'
$ws.Range("E19").Value = '  +3.20%  '
$ws.Range("D20").Value = '''8.03'
$ws.Range("E20").Value = '  +2.11%  '
$ws.Range("D21").Value = '''372.21'
$ws.Range("E21").Value = '  -2.05%  '
$ws.Range("D22").Value = '''4.19'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").Value = '''2.06'
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").Value = '''71.88'
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  -2.67%  '
$ws.Range("D27").Value = '''9.69'
$ws.Range("E27").Value = '  -4.19%  '
$ws.Range("D28").Value = '2.795.64'
$ws.Range("E28").Value = '  +1.06%  '
$ws.Range("D29").Value = '''0.999'
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").Value = '''8.15'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '''497.71'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("E33").Value = '  -2.21%  '
$ws.Range("D34").Value = '''1.82'
$ws.Range("E34").Value = '  -0.63%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '''162.90'
$ws.Range("E36").Value = '  -1.03%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.117'
$ws.Range("E37").Value = '  +2.95%  '
$ws.Range("B38").Value = 'EthereumClassic'
$ws.Range("C38").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D38").Value = '''19.44'
$ws.Range("E38").Value = '  +0.72%  '
$ws.Range("D39").Value = '''18.92'
$ws.Range("E39").Value = '  -0.85%  '
$ws.Range("E40").Value = '  -2.56%  '
$ws.Range("E42").Value = '  -6.04%  '
$ws.Range("D43").Value = '''2.57'
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("D44").Value = '''4.89'
$ws.Range("E44").Value = '  -3.52%  '
$ws.Range("D45").Value = '''0.330'
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").Value = '''39.19'
$ws.Range("E46").Value = '  -0.33%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '''154.73'
$ws.Range("E47").Value = '  +2.62%  '
$ws.Range("E48").Value = '  +0.06%  '
$ws.Range("D49").Value = '''0.552'
$ws.Range("E49").Value = '  +1.52%  '
$ws.Range("E50").Value = '  +0.58%  '
$ws.Range("D51").Value = '''0.0752'
$ws.Range("E51").Value = '  -1.28%  '
